# chore: update Sheets via scheduled runner
#
# Refresh cached market-board snapshot values (currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW
# and WVR sheets, as produced by the scheduled data-refresh job.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 15
$ws.Range("H15").Value = 2458.4746
$ws.Range("I15").Value = 2458.4746
$ws.Range("K15").Value = 7375.4238
$ws.Range("M15").Value = -7206.4238

# Row 46
$ws.Range("H46").Value = 91977.45
$ws.Range("I46").Value = 744.6667
$ws.Range("J46").Value = 126189.75
$ws.Range("K46").Value = 2234.0001
$ws.Range("L46").Value = 378569.25
$ws.Range("M46").Value = -2115.0001
$ws.Range("N46").Value = -378807.25

# Row 60
$ws.Range("H60").Value = 91977.45
$ws.Range("I60").Value = 744.6667
$ws.Range("J60").Value = 126189.75
$ws.Range("K60").Value = 2234.0001
$ws.Range("L60").Value = 378569.25
$ws.Range("M60").Value = -1750.0001
$ws.Range("N60").Value = -379537.25

# Row 76
$ws.Range("H76").Value = 3124.52
$ws.Range("I76").Value = 2793.0476
$ws.Range("J76").Value = 4864.75
$ws.Range("K76").Value = 2793.0476
$ws.Range("L76").Value = 4864.75
$ws.Range("M76").Value = -2478.0476
$ws.Range("N76").Value = -5494.75

# Row 79
$ws.Range("H79").Value = 3124.52
$ws.Range("I79").Value = 2793.0476
$ws.Range("J79").Value = 4864.75
$ws.Range("K79").Value = 2793.0476
$ws.Range("L79").Value = 4864.75
$ws.Range("M79").Value = -1701.0476
$ws.Range("N79").Value = -7048.75

# Row 132
$ws.Range("H132").Value = 3032373.5
$ws.Range("I132").Value = 3391690.8
$ws.Range("J132").Value = 3842.5715
$ws.Range("K132").Value = 10175072.4
$ws.Range("L132").Value = 11527.7145
$ws.Range("M132").Value = -10172542.4
$ws.Range("N132").Value = -16587.7145

# Row 135
$ws.Range("H135").Value = 1655.6316
$ws.Range("I135").Value = 1341.0625
$ws.Range("J135").Value = 3333.3333
$ws.Range("K135").Value = 12069.5625
$ws.Range("L135").Value = 29999.9997
$ws.Range("M135").Value = -9534.5625
$ws.Range("N135").Value = -35069.9997

# Row 137
$ws.Range("H137").Value = 3229004.5
$ws.Range("I137").Value = 4765340.5
$ws.Range("J137").Value = 2699
$ws.Range("K137").Value = 14296021.5
$ws.Range("L137").Value = 8097
$ws.Range("M137").Value = -14293471.5
$ws.Range("N137").Value = -13197

# Row 138
$ws.Range("H138").Value = 2512.2222
$ws.Range("I138").Value = 901.2826
$ws.Range("J138").Value = 3910.3962
$ws.Range("K138").Value = 2703.8478
$ws.Range("L138").Value = 11731.1886
$ws.Range("M138").Value = 2436.1522
$ws.Range("N138").Value = -22011.1886

# Row 141
$ws.Range("H141").Value = 294111.97
$ws.Range("I141").Value = 1376.862
$ws.Range("J141").Value = 1237369.5
$ws.Range("K141").Value = 4130.586
$ws.Range("L141").Value = 3712108.5
$ws.Range("M141").Value = 1049.414
$ws.Range("N141").Value = -3722468.5


$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 3657.75
$ws.Range("I32").Value = 3060.7805
$ws.Range("J32").Value = 6377.278
$ws.Range("K32").Value = 3060.7805
$ws.Range("L32").Value = 6377.278
$ws.Range("M32").Value = -2773.7805
$ws.Range("N32").Value = -6951.278

# Row 61
$ws.Range("H61").Value = 1476.94
$ws.Range("I61").Value = 639.4878
$ws.Range("J61").Value = 5292
$ws.Range("K61").Value = 639.4878
$ws.Range("L61").Value = 5292
$ws.Range("M61").Value = -427.4878
$ws.Range("N61").Value = -5716

# Row 74
$ws.Range("H74").Value = 701.64703
$ws.Range("I74").Value = 701.64703
$ws.Range("K74").Value = 701.64703
$ws.Range("M74").Value = 172.35297

# Row 77
$ws.Range("H77").Value = 701.64703
$ws.Range("I77").Value = 701.64703
$ws.Range("K77").Value = 3508.23515
$ws.Range("M77").Value = 859.76485

# Row 102
$ws.Range("H102").Value = 5689.8335
$ws.Range("I102").Value = 6121
$ws.Range("J102").Value = 4827.5
$ws.Range("K102").Value = 6121
$ws.Range("L102").Value = 4827.5
$ws.Range("M102").Value = -4499
$ws.Range("N102").Value = -8071.5

# Row 132
$ws.Range("H132").Value = 2065.1
$ws.Range("I132").Value = 1549.7693
$ws.Range("K132").Value = 4649.3079
$ws.Range("M132").Value = -2119.3079

# Row 136
$ws.Range("H136").Value = 1476.94
$ws.Range("I136").Value = 639.4878
$ws.Range("J136").Value = 5292
$ws.Range("K136").Value = 1918.4634
$ws.Range("L136").Value = 15876
$ws.Range("M136").Value = 631.5365999999999
$ws.Range("N136").Value = -20976


$ws = $wb.Worksheets.Item("BSM")

# Row 7
$ws.Range("H7").Value = 16614
$ws.Range("I7").Value = 1500
$ws.Range("J7").Value = 22659.6
$ws.Range("K7").Value = 1500
$ws.Range("L7").Value = 22659.6
$ws.Range("M7").Value = -1387
$ws.Range("N7").Value = -22885.6

# Row 38
$ws.Range("H38").Value = 70036
$ws.Range("J38").Value = 70036
$ws.Range("L38").Value = 70036
$ws.Range("N38").Value = -70868

# Row 105
$ws.Range("H105").Value = 1597.3334
$ws.Range("I105").Value = 1430.7142
$ws.Range("J105").Value = 2180.5
$ws.Range("K105").Value = 1430.7142
$ws.Range("L105").Value = 2180.5
$ws.Range("M105").Value = 316.2858000000001
$ws.Range("N105").Value = -5674.5

# Row 107
$ws.Range("H107").Value = 4049.4211
$ws.Range("I107").Value = 3339
$ws.Range("J107").Value = 4838.778
$ws.Range("K107").Value = 3339
$ws.Range("L107").Value = 4838.778
$ws.Range("M107").Value = -1419
$ws.Range("N107").Value = -8678.778

# Row 132
$ws.Range("H132").Value = 40000
$ws.Range("J132").Value = 40000
$ws.Range("L132").Value = 40000
$ws.Range("N132").Value = -50120

# Row 134
$ws.Range("H134").Value = 1990.5264
$ws.Range("I134").Value = 1524.2554
$ws.Range("J134").Value = 4182
$ws.Range("K134").Value = 4572.7662
$ws.Range("L134").Value = 12546
$ws.Range("M134").Value = -2037.7662
$ws.Range("N134").Value = -17616


$ws = $wb.Worksheets.Item("CRP")

# Row 10
$ws.Range("H10").Value = 7146.7915
$ws.Range("I10").Value = 540.4666999999999
$ws.Range("J10").Value = 18157.334
$ws.Range("K10").Value = 540.4666999999999
$ws.Range("L10").Value = 18157.334
$ws.Range("M10").Value = -401.4666999999999
$ws.Range("N10").Value = -18435.334

# Row 31
$ws.Range("H31").Value = 2860423.5
$ws.Range("I31").Value = 4763617
$ws.Range("J31").Value = 5632.9287
$ws.Range("K31").Value = 4763617
$ws.Range("L31").Value = 5632.9287
$ws.Range("M31").Value = -4763322
$ws.Range("N31").Value = -6222.9287

# Row 34
$ws.Range("H34").Value = 2860423.5
$ws.Range("I34").Value = 4763617
$ws.Range("J34").Value = 5632.9287
$ws.Range("K34").Value = 4763617
$ws.Range("L34").Value = 5632.9287
$ws.Range("M34").Value = -4763415
$ws.Range("N34").Value = -6036.9287

# Row 58
$ws.Range("H58").Value = 9261869
$ws.Range("I58").Value = 1590.125
$ws.Range("J58").Value = 35719810
$ws.Range("K58").Value = 1590.125
$ws.Range("L58").Value = 35719810
$ws.Range("M58").Value = -1387.125
$ws.Range("N58").Value = -35720216

# Row 122
$ws.Range("H122").Value = 3608.28
$ws.Range("I122").Value = 2703.25
$ws.Range("K122").Value = 8109.75
$ws.Range("M122").Value = -5659.75

# Row 132
$ws.Range("H132").Value = 1668.4894
$ws.Range("I132").Value = 1223.55
$ws.Range("J132").Value = 4211
$ws.Range("K132").Value = 3670.65
$ws.Range("L132").Value = 12633
$ws.Range("M132").Value = -1140.65
$ws.Range("N132").Value = -17693

# Row 134
$ws.Range("H134").Value = 1503.262
$ws.Range("I134").Value = 828.6177
$ws.Range("J134").Value = 4370.5
$ws.Range("K134").Value = 2485.8531
$ws.Range("L134").Value = 13111.5
$ws.Range("M134").Value = 49.14689999999973
$ws.Range("N134").Value = -18181.5

# Row 136
$ws.Range("H136").Value = 9261869
$ws.Range("I136").Value = 1590.125
$ws.Range("J136").Value = 35719810
$ws.Range("K136").Value = 4770.375
$ws.Range("L136").Value = 107159430
$ws.Range("M136").Value = -2220.375
$ws.Range("N136").Value = -107164530


$ws = $wb.Worksheets.Item("CUL")

# Row 68
$ws.Range("H68").Value = 4100.222
$ws.Range("I68").Value = 699.5
$ws.Range("J68").Value = 5071.857
$ws.Range("K68").Value = 2098.5
$ws.Range("L68").Value = 15215.571
$ws.Range("M68").Value = -1287.5
$ws.Range("N68").Value = -16837.571

# Row 71
$ws.Range("H71").Value = 4100.222
$ws.Range("I71").Value = 699.5
$ws.Range("J71").Value = 5071.857
$ws.Range("K71").Value = 6295.5
$ws.Range("L71").Value = 45646.713
$ws.Range("M71").Value = -2239.5
$ws.Range("N71").Value = -53758.713

# Row 117
$ws.Range("H117").Value = 2119.8
$ws.Range("I117").Value = 211
$ws.Range("J117").Value = 2331.889
$ws.Range("K117").Value = 633
$ws.Range("L117").Value = 6995.667
$ws.Range("M117").Value = 2809
$ws.Range("N117").Value = -13879.667

# Row 132
$ws.Range("H132").Value = 3517.3333
$ws.Range("J132").Value = 5000
$ws.Range("L132").Value = 45000
$ws.Range("N132").Value = -50060


$ws = $wb.Worksheets.Item("GSM")

# Row 4
$ws.Range("H4").Value = 49223.332
$ws.Range("I4").Value = 2000
$ws.Range("K4").Value = 2000
$ws.Range("M4").Value = -1888

# Row 58
$ws.Range("H58").Value = 30015.334
$ws.Range("J58").Value = 30015.334
$ws.Range("L58").Value = 30015.334
$ws.Range("N58").Value = -30569.334

# Row 132
$ws.Range("H132").Value = 1932.638
$ws.Range("I132").Value = 1486.2222
$ws.Range("J132").Value = 3477.923
$ws.Range("K132").Value = 4458.6666
$ws.Range("L132").Value = 10433.769
$ws.Range("M132").Value = -1928.6666
$ws.Range("N132").Value = -15493.769


$ws = $wb.Worksheets.Item("LTW")

# Row 132
$ws.Range("H132").Value = 2132.2703
$ws.Range("I132").Value = 1411.76
$ws.Range("K132").Value = 4235.28
$ws.Range("M132").Value = -1705.28

# Row 135
$ws.Range("H135").Value = 29796.072
$ws.Range("J135").Value = 29796.072
$ws.Range("L135").Value = 29796.072
$ws.Range("N135").Value = -39936.072

# Row 136
$ws.Range("H136").Value = 2440944.2
$ws.Range("I136").Value = 2942624
$ws.Range("J136").Value = 4214.2856
$ws.Range("K136").Value = 8827872
$ws.Range("L136").Value = 12642.8568
$ws.Range("M136").Value = -8825322
$ws.Range("N136").Value = -17742.8568


$ws = $wb.Worksheets.Item("WVR")

# Row 92
$ws.Range("H92").Value = 39000
$ws.Range("J92").Value = 39000
$ws.Range("L92").Value = 39000
$ws.Range("N92").Value = -43992

# Row 100
$ws.Range("H100").Value = 1453.7059
$ws.Range("I100").Value = 1171.4286
$ws.Range("J100").Value = 1651.3
$ws.Range("K100").Value = 2342.8572
$ws.Range("L100").Value = 3302.6
$ws.Range("M100").Value = -1801.8572
$ws.Range("N100").Value = -4384.6

# Row 123
$ws.Range("H123").Value = 15000
$ws.Range("J123").Value = 15000
$ws.Range("L123").Value = 15000
$ws.Range("N123").Value = -24800

# Row 132
$ws.Range("H132").Value = 216691.88
$ws.Range("I132").Value = 265497.44
$ws.Range("J132").Value = 48090.816
$ws.Range("K132").Value = 796492.3200000001
$ws.Range("L132").Value = 144272.448
$ws.Range("M132").Value = -793962.3200000001
$ws.Range("N132").Value = -149332.448

# Row 136
$ws.Range("H136").Value = 1521.68
$ws.Range("I136").Value = 1017
$ws.Range("J136").Value = 1918.2142
$ws.Range("K136").Value = 3051
$ws.Range("L136").Value = 5754.642599999999
$ws.Range("M136").Value = -501
$ws.Range("N136").Value = -10854.6426

